# Applies the diff: reorders the "HOP 4.7"/"HOP 4.8" caption paragraphs so
# each sits *before* its screenshot (instead of after), adds two new blank
# paragraphs right after "HOP 4.6", stamps a <w:lastRenderedPageBreak/> on
# the Picture 4 run, and refreshes the wp14:anchorId/editId pairs on the
# Picture 4/5/7 drawings.

$d = $word.ActiveDocument

$xmlns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $xmlns + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Picture 7 (just above "HOP 4.6"): refresh anchorId/editId only.
# ---------------------------------------------------------------------
$pic7Body = '<w:p><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="6B36240B" wp14:editId="1711B901"><wp:extent cx="2634799" cy="2794715"/><wp:effectExtent l="0" t="0" r="0" b="5715"/><wp:docPr id="7" name="Picture 7" descr="A picture containing text, screenshot, display, electronics&#10;&#10;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="7" name="Picture 7" descr="A picture containing text, screenshot, display, electronics&#10;&#10;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2642194" cy="2802558"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$pic7Para = $d.Paragraphs.Item(12)
$pic7Range = $d.Range($pic7Para.Range.Start, $pic7Para.Range.End)
$pic7Range.InsertXML((New-PkgXml $pic7Body))

# ---------------------------------------------------------------------
# Paragraph layout is unchanged up through "HOP 4.6" (index 13). What
# follows today is: [Picture 4] [blank] [HOP 4.7] [Picture 5] [HOP 4.8]
# [blank] [HOP 4.9]. Work from the bottom up so earlier indices stay valid.
# ---------------------------------------------------------------------

# 2) Drop the trailing blank paragraph after "HOP 4.8" (index 19).
$d.Paragraphs.Item(19).Range.Delete()

# 3) Drop the old "HOP 4.8" paragraph (index 18); its text gets rebuilt
#    right before Picture 5 in step 5.
$d.Paragraphs.Item(18).Range.Delete()

# 4) Drop the old "HOP 4.7" paragraph (index 16); rebuilt before Picture 4
#    in step 6.
$d.Paragraphs.Item(16).Range.Delete()

# 5) Drop the blank paragraph that used to sit between Picture 4 and
#    "HOP 4.7" (index 15).
$d.Paragraphs.Item(15).Range.Delete()

# Paragraphs 13-16 are now: "HOP 4.6", [Picture 4], [Picture 5], "HOP 4.9".

# ---------------------------------------------------------------------
# 6) Replace Picture 5 (now index 15) with itself + refreshed anchor ids,
#    preceded by a brand-new "HOP 4.8" paragraph.
# ---------------------------------------------------------------------
$hop48Body = '<w:p><w:r><w:t xml:space="preserve">HOP </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>.8</w:t></w:r></w:p>'
$pic5Body = '<w:p><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="4FE09F53" wp14:editId="774D1634"><wp:extent cx="5206405" cy="2987898"/><wp:effectExtent l="0" t="0" r="0" b="3175"/><wp:docPr id="5" name="Picture 5" descr="Graphical user interface, website&#10;&#10;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="5" name="Picture 5" descr="Graphical user interface, website&#10;&#10;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5211396" cy="2990762"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$pic5Para = $d.Paragraphs.Item(15)
$pic5Range = $d.Range($pic5Para.Range.Start, $pic5Para.Range.End)
$pic5Range.InsertXML((New-PkgXml ($hop48Body + $pic5Body)))

# ---------------------------------------------------------------------
# 7) Replace Picture 4 (index 14) with itself + refreshed anchor ids and
#    a new <w:lastRenderedPageBreak/>, preceded by two new blank
#    paragraphs and a brand-new "HOP 4.7" paragraph.
# ---------------------------------------------------------------------
$blanksBody = '<w:p/><w:p/>'
$hop47Body = '<w:p><w:r><w:t xml:space="preserve">HOP </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>.7</w:t></w:r></w:p>'
$pic4Body = '<w:p><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="528A4ED4" wp14:editId="1EFF0678"><wp:extent cx="4185634" cy="2500916"/><wp:effectExtent l="0" t="0" r="5715" b="0"/><wp:docPr id="4" name="Picture 4" descr="Graphical user interface, text, website&#10;&#10;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="4" name="Picture 4" descr="Graphical user interface, text, website&#10;&#10;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId9"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="4188096" cy="2502387"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$pic4Para = $d.Paragraphs.Item(14)
$pic4Range = $d.Range($pic4Para.Range.Start, $pic4Para.Range.End)
$pic4Range.InsertXML((New-PkgXml ($blanksBody + $hop47Body + $pic4Body)))

Write-Output "done"
